$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at row 60 (pushes the old row 60 -> row 61), and grow the
# table (and its autofilter range) by one row so the table keeps covering
# the whole data range, just like Excel does when a row is inserted inside
# a table.
$ws.Rows("60:60").Insert()
$lo.Resize($ws.Range("A1:E61"))

# Populate the new row with the Word.Style enum snippet mapping.
$ws.Range("A60").Value = "Style"
$ws.Range("C60").Value = "enum"
$ws.Range("D60").Value = "word-paragraph-insert-formatted-text"
$ws.Range("E60").Value = "addPreStyledFormattedText"

# Match the cursor/selection position left behind by the edit.
$ws.Range("D61").Select()
